$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '87.739.14'
$ws.Range("E2").Value = '  +3.65%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.264.14'
$ws.Range("E3").Value = '  -1.06%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.99'
$ws.Range("E5").Value = '  -3.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '627.78'
$ws.Range("E6").Value = '  -1.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.378'
$ws.Range("E7").Value = '  +16.81%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.692'
$ws.Range("E8").Value = '  +16.69%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.258.33'
$ws.Range("E10").Value = '  -1.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.578'
$ws.Range("E11").Value = '  -2.76%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.183'
$ws.Range("E12").Value = '  +9.78%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000259'
$ws.Range("E13").Value = '  -7.25%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.15'
$ws.Range("E14").Value = '  -0.36%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.858.04'
$ws.Range("E15").Value = '  -0.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.33'
$ws.Range("E16").Value = '  -1.86%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.281.24'
$ws.Range("E17").Value = '  +3.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.256.93'
$ws.Range("E18").Value = '  -0.85%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.13'
$ws.Range("E19").Value = '  -2.54%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.05'
$ws.Range("E20").Value = '  -3.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '434.21'
$ws.Range("E21").Value = '  -0.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.90'
$ws.Range("E22").Value = '  -3.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.33'
$ws.Range("E23").Value = '  +2.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.28'
$ws.Range("E24").Value = '  -2.46%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.48'
$ws.Range("E25").Value = '  +2.40%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.13'
$ws.Range("E26").Value = '  -7.35%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.380.58'
$ws.Range("E27").Value = '  -1.90%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '76.53'
$ws.Range("E28").Value = '  -2.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000130'
$ws.Range("E29").Value = '  -1.45%  '

$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.179'
$ws.Range("E31").Value = '  +9.75%  '

$ws.Range("E32").Value = '  +0.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.77'
$ws.Range("E33").Value = '  -5.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '548.76'
$ws.Range("E34").Value = '  -8.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.39'
$ws.Range("E35").Value = '  -11.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.96'
$ws.Range("E36").Value = '  -3.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.97'
$ws.Range("E37").Value = '  +8.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.137'
$ws.Range("E38").Value = '  -10.68%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.48'
$ws.Range("E39").Value = '  -3.36%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.09%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '21.72'
$ws.Range("E41").Value = '  +3.64%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.01'
$ws.Range("E42").Value = '  -2.53%  '

$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.393'
$ws.Range("E43").Value = '  -5.55%  '

$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.94'
$ws.Range("E44").Value = '  -5.35%  '

$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '154.97'
$ws.Range("E46").Value = '  -2.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '179.61'
$ws.Range("E47").Value = '  -5.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.90'
$ws.Range("E48").Value = '  -0.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.29'
$ws.Range("E49").Value = '  -4.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.22'
$ws.Range("E50").Value = '  -0.73%  '

$ws.Range("E51").Value = '  +10.73%  '
